$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the affected rows
$ws.Range("F3").Value = 1
$ws.Range("F7").Value = -5
$ws.Range("F11").Value = -1
$ws.Range("F21").Value = -4
$ws.Range("F26").Value = -9
$ws.Range("F29").Value = 8
$ws.Range("F33").Value = -4
$ws.Range("F34").Value = -4
